$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above the current row 345, shifting the existing
# rows 345:375 down to 347:377 (dimension grows from T375 to T377).
$ws.Rows("345:346").Insert()

# Populate the two newly inserted rows with the new weekly price records.
$ws.Range("A345").Value = 8
$ws.Range("B345").Value = "Terminal La Palmera de La Serena"
$ws.Range("C345").Value = "Coquimbo"
$ws.Range("D345").Value = 44449
$ws.Range("E345").Value = 4
$ws.Range("F345").Value = "Fruta"
$ws.Range("G345").Value = 100101
$ws.Range("H345").Value = "Berries"
$ws.Range("I345").Value = 100112025
$ws.Range("J345").Value = "Frutilla"
$ws.Range("K345").Value = "Sin especificar"
$ws.Range("L345").Value = "Primera"
$ws.Range("M345").Value = 300
$ws.Range("N345").Value = 26000
$ws.Range("O345").Value = 27000
$ws.Range("P345").Value = 26500
$ws.Range("Q345").Value = "$/bandeja 7 kilos"
$ws.Range("R345").Value = "Provincia de Melipilla"
$ws.Range("S345").Value = 3786
$ws.Range("T345").Value = 7

$ws.Range("A346").Value = 8
$ws.Range("B346").Value = "Terminal La Palmera de La Serena"
$ws.Range("C346").Value = "Coquimbo"
$ws.Range("D346").Value = 44449
$ws.Range("E346").Value = 4
$ws.Range("F346").Value = "Fruta"
$ws.Range("G346").Value = 100101
$ws.Range("H346").Value = "Berries"
$ws.Range("I346").Value = 100112025
$ws.Range("J346").Value = "Frutilla"
$ws.Range("K346").Value = "Sin especificar"
$ws.Range("L346").Value = "Segunda"
$ws.Range("M346").Value = 200
$ws.Range("N346").Value = 20000
$ws.Range("O346").Value = 21000
$ws.Range("P346").Value = 20500
$ws.Range("Q346").Value = "$/bandeja 7 kilos"
$ws.Range("R346").Value = "Provincia de Melipilla"
$ws.Range("S346").Value = 2929
$ws.Range("T346").Value = 7
